# Removal of the use of the LPS25 sensor
$wb = $excel.ActiveWorkbook

$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# Typography sheet: "Default" row no longer uses the LPS25-specific 30px size
# nor the wildcard/fallback customizations that supported it.
$wsTypo.Range("D4").Value = 20

$wsTypo.Range("G4:J6").ClearContents()
$wsTypo.Range("G4:J6").Style = "Normal"

# Translation sheet: remove the two translation rows (LPS25 temperature and
# pressure readouts) that referenced the sensor.
$wsTrans.Rows.Item(5).Delete()
$wsTrans.Rows.Item(4).Delete()
